$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '26.877.95'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.98%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.872.42'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -1.45%  '
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '301.92'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.66%  '
$ws.Range('E6').Value = '  -0.21%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5299'
$ws.Range('D7').Style = "Normal"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3752'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -1.31%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07157'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -1.76%  '
$ws.Range('E10').Value = '  +1.20%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.8845'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -2.37%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.08137'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -0.95%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.835.06'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -3.24%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '92.96'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -2.68%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '5.270'
$ws.Range('D15').Style = "Normal"
$ws.Range('E16').Value = '  -0.07%  '
$ws.Range('E17').Value = '  +0.23%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000008530'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -1.45%  '
$ws.Range('E19').Value = '  -0.24%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '26.892.30'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -1.08%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '4.970'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -2.92%  '
$ws.Range('E22').Value = '  -1.11%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '6.364'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -1.65%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '147.30'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -1.53%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.260'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -2.92%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.735'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.43%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '18.00'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -1.49%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '114.42'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.85%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '4.736'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -1.81%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.567'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -6.28%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.09097'
$ws.Range('D31').Style = "Normal"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.7967'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.35%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.04976'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -1.40%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.988'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +0.84%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.172'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -4.31%  '
$ws.Range('B36').Value = 'TheSandbox'
$ws.Range('C36').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.5845'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +1.87%  '
$ws.Range('B37').Value = 'MXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.197'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -5.41%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.601'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.55%  '
$ws.Range('E39').Value = '  -0.88%  '
$ws.Range('E40').Value = '  -2.23%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '6.592'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.49%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '8.875'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -1.82%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '115.64'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.65%  '
$ws.Range('E44').Value = '  +2.96%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.1491'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -1.71%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.9997'
$ws.Range('D46').Style = "Normal"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '9.976'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -1.38%  '
$ws.Range('E48').Value = '  -1.91%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '37.96'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -1.64%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.06029'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +1.23%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '62.31'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -2.86%  '
